# Add a new worksheet "newSheet" after the existing "WebTable" sheet and
# populate its header row (A1:R1) with the Company/contact/country labels
# (reusing the same shared strings already used on the WebTable sheet),
# giving each header cell the same highlighted fill used elsewhere in the
# workbook. Finish by reselecting the WebTable sheet/cell so it remains the
# active tab, matching the final state of the workbook.

$wb = $excel.ActiveWorkbook
$webTable = $wb.Worksheets.Item("WebTable")

# Keep WebTable active while we create the new sheet so it is inserted
# immediately after it (i.e. becomes the 3rd / last tab).
$webTable.Activate()
$ws = $wb.Worksheets.Add($null, $webTable)
$ws.Name = "newSheet"

$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R")
$headers = @("Company", "Company", "Company", "Company", "Company", "Company", `
             "contact", "contact", "contact", "contact", "contact", "contact", `
             "country", "country", "country", "country", "country", "country")

for ($i = 0; $i -lt $columns.Length; $i++) {
    $cell = $ws.Range($columns[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Interior.ColorIndex = 55
}

# Return focus to the WebTable sheet / cell, which remains the active tab.
$webTable.Activate()
$webTable.Range("H17").Select()
